$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 582 (pushes the existing 582:607 block down to
# 583:608, matching Excel's default "insert" behaviour of copying the format
# of the row above into the freshly inserted blank row).
$ws.Rows.Item(582).Insert()

# Populate the newly inserted row 582 with the new weekly price-report entry.
$ws.Cells.Item(582, 1).Value = 6
$ws.Cells.Item(582, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(582, 3).Value = "Metropolitana"
$ws.Cells.Item(582, 4).Value = 45147
$ws.Cells.Item(582, 5).Value = 13
$ws.Cells.Item(582, 6).Value = 100112032
$ws.Cells.Item(582, 7).Value = "Zapallo italiano"
$ws.Cells.Item(582, 8).Value = "Sin especificar"
$ws.Cells.Item(582, 9).Value = "Primera"
$ws.Cells.Item(582, 10).Value = 420
$ws.Cells.Item(582, 11).Value = 15000
$ws.Cells.Item(582, 12).Value = 16000
$ws.Cells.Item(582, 13).Value = 15238
$ws.Cells.Item(582, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(582, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(582, 16).Value = 305
$ws.Cells.Item(582, 17).Value = 50
$ws.Cells.Item(582, 18).Value = "Hortaliza"
